$wb = $excel.ActiveWorkbook

# --- Sheet "Oct-24 RT Pk" ---
$ws1 = $wb.Worksheets.Item("Oct-24 RT Pk")
$ws1.Range("O5").Value = 'Severe, unprecedented derate. Clean relationship between binding events and the derate. Hasn''t bound since.'
$ws1.Range("O6").Value = 'Mt. Vernon-W. Frankfort 230 kV outage. Moderate sum of wind generation and coal generation. Load to the South, Southeast.'
$ws1.Range("O7").Value = 'Bound often before 12/17/2024. A clear regime change occurred on 12/18/2024 which drastically reduced post-contingent flow.'
$ws1.Range("O8").Value = 'Diamond Trail-Hills 345 kV, Powesheik Co-Reasnor 161 kV outages. Wind push.'
$ws1.Range("O9").Value = 'Big Stone TR 230/115 kV, Benson-Swenoda 115 kV outages increased binding likelihood.'
$ws1.Range("O10").Value = 'Weak wind in OTP and high Boswell generation.'
$ws1.Range("O11").Value = 'Doud-Boone Junction 161 kV, Karma-Perry 161 kV outages. Story County, Great Pathfinder wind generation.'
$ws1.Range("O12").Value = 'Big Stone TR 230/115 kV, Big Stone-Brownsville 230 kV outages. High sensitivity to wind push at Dakota Range, Crowned Ridge, Deuel Harvest, etc. Clean relationship with Big Stone TR outage.'
$ws1.Range("O13").Value = 'Sarepta-Longwood 345 kV outage. ERCOT East DC tie line "load."'
$ws1.Range("O14").Value = 'Doud-Boone Junction 161 kV, Karma-Perry 161 kV outages. Story County, Great Pathfinder wind generation.'
$ws1.Range("O15").Value = 'Monticello-Lafayette-Lafayette LC 138 kV outages.'
$ws1.Range("P15").Value = 'Bound briefly before substantial uprate on 10/29/2025 from 198 MW to 357 MW (post-contingent). Unlikely to bind in the current regime.'
$ws1.Range("O16").Value = 'Cordova-Substation 39 345 kV outage. Quad Cities push.'
$ws1.Range("O17").Value = 'Sarepta-Longwood 345 kV outage. Weak Harrison County generation. Weak OK wind generation.'
$ws1.Range("O18").Value = 'Harmony-Lansing 161 kV outage.'
$ws1.Range("P18").Value = 'This is odd, as it looks like the monitored element is radial. I''m not fully sure what''s going on here.'
$ws1.Range("O19").Value = 'Floyd-Emery-Sheffield, Worth Co-Glenville-Hayward 161 kV outages. Wind generation in Northern IA.'
$ws1.Range("O20").Value = 'Southward wind transfer from IA into MO. Hills-Diamond Trails 345 kV outage.'
$ws1.Range("O21").Value = 'Eastward wind transfer from South Central IA to Southeast IA and IL. Hills-Diamond Trail 345 kV outage.  High wind generation. Excess MEC, ALTW generation.'

# --- Sheet "Oct-24 RT Off" ---
$ws2 = $wb.Worksheets.Item("Oct-24 RT Off")
$ws2.Range("O3").Value = 'Mt. Vernon-W. Frankfort 230 kV outage. Moderate sum of wind generation and coal generation. Load to the South, Southeast.'
$ws2.Range("O5").Value = 'Severe, unprecedented derate. Clean relationship between binding events and the derate. Hasn''t bound since.'
$ws2.Range("O6").Value = 'Diamond Trail-Hills 345 kV, Powesheik Co-Reasnor 161 kV outages. Wind push.'
$ws2.Range("O8").Value = 'Sarepta-Longwood 345 kV outage. Weak Harrison County generation. Weak OK wind generation.'
$ws2.Range("O9").Value = 'Thomas Hill-Mead 161 kV outage. Heavy SPP wind generation. Thomas Hill outage.'
$ws2.Range("O10").Value = 'White Cloud, Clear Creek, Rock Creek wind generation. Southwest IA wind generation. Excess load in St. Joseph, Kansas City.'
$ws2.Range("O11").Value = 'Excess ComEd generation. Excess MEC generation. Sheffield-Gary Ave-Dune Acres 345 kV outages. Chicago Ave-Gary Ave 138 kV outage.'
$ws2.Range("O12").Value = 'Southward wind transfer from IA into MO. Hills-Diamond Trails 345 kV outage.'
$ws2.Range("O13").Value = 'Reasnor-Poweshiek Co-Irvine 161 kV outages. Hills-Diamond Trail 345 kV outage. Knoxville-Lucas 69 kV outages. Grinnell-S. Grinnell 69 kV outage.'
$ws2.Range("O14").Value = 'Southward wind transfer from IA into MO. Hills-Diamond Trails 345 kV outage. High wind generation. Low Ottumwa output.'
$ws2.Range("O15").Value = 'Big Stone TR 230/115 kV, Big Stone-Brownsville 230 kV outages. High sensitivity to wind push at Dakota Range, Crowned Ridge, Deuel Harvest, etc. Clean relationship with Big Stone TR outage.'
$ws2.Range("O16").Value = 'Winger-Riverton, Hubbard-Audubon 230 kV outages. Wind generation to the West.'
$ws2.Range("O17").Value = 'Eastward wind transfer from South Central IA to Southeast IA and IL. Hills-Diamond Trail 345 kV outage.  High wind generation. Excess MEC, ALTW generation.'
$ws2.Range("O18").Value = 'Clear relationship with Hazleton-Arnold 345 kV outage, which bottlenecks wind supply from IA and forces it over the constraint.'
$ws2.Range("O19").Value = 'Maryville-Maryville 161 kV outage (extreme impact). Highly sensitive to Clear Creek, White Cloud wind push in the post-contingent setting with Maryville-Maryville 161 kV OOS.'
$ws2.Range("O20").Value = 'Palisades-Argenta-Tompkins, Battle Ck-Oneida 345 kV outages. Verona-Convis 138 kV outage.'
$ws2.Range("P20").Value = 'Will not bind when Calhoun Solar, Cereal City Solar are generating.'
$ws2.Range("O21").Value = 'Wind transfer stepping down to 115 kV level. Clean relationship with Chub Lk-Hampton 345 kV outage. Excess MEC supply.'

# --- Restore selections / active cell per final state ---
$ws2.Range("O22").Select()
$ws1.Select()
$ws1.Range("N30").Select()
